# Auto-generated Excel COM-interop edit script
# Updates cryptos price/volume values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.386.34'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.565.78'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D5").Value = '210.90'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '44.56'
$ws.Range("E8").Value = '  -3.38%  '
$ws.Range("D9").Value = '23.67'
$ws.Range("E9").Value = '  -1.45%  '
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = '0.0589'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '0.0894'
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").Value = '1.789.78'
$ws.Range("D14").Value = '1.560.67'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '3.67'
$ws.Range("E15").Value = '  -0.16%  '
$ws.Range("D16").Value = '28.360.33'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '0.514'
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").Value = '60.73'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '228.47'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = '0.0' + [char]8323 + '0681'
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '3.95'
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").Value = '8.95'
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").Value = '2.06'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").Value = '150.47'
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("D27").Value = '14.90'
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").Value = '0.0477'
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("E32").Value = '  -4.22%  '
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("D34").Value = '3.09'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '1.387.12'
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  -3.05%  '
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("D41").Value = '0.521'
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '0.786'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").Value = '5.35'
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("D47").Value = '62.27'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").Value = '0.917'
$ws.Range("E48").Value = '  -6.01%  '
$ws.Range("D49").Value = '1.702.18'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").Value = '85.56'
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("E51").Value = '  -1.88%  '
